# James's Week 6 and 7 Timesheets
# Roll the timesheet forward to the week commencing 08/02/2014 and fill in
# the hours worked, plus mark the employee signature date for the next
# (week 7) sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week-of date (G8) - Saturday 08/02/2014, the days below (A11:A17) are
# derived from this via formula. Use the serial date number (41678 ==
# 08/02/2014) since the cell is already formatted as a date (numFmtId 14).
$ws.Range("G8").Value = 41678

# Regular hours worked each day that week.
$ws.Range("B11").Value = 3
$ws.Range("B12").Value = 4
$ws.Range("B13").Value = 1
$ws.Range("B14").Value = 2
$ws.Range("B15").Value = 0
$ws.Range("B16").Value = 3
$ws.Range("B17").Value = 2

# Employee signed off - date next due (week 7 sheet date).
$ws.Range("D27").Value = "28/02/2014"

# Clear the old supervisor-signature date, to be re-signed for this week.
$ws.Range("D29").ClearContents()

# Select the next signature-date cell, matching where the author left off
# editing (bottom section of the form).
$ws.Range("D28:E28").Select() | Out-Null
